# Set line spacing to single (1.0) for every paragraph in the document,
# matching Word's "Line Spacing: Single" formatting option
# (OOXML: <w:spacing w:line="240" w:lineRule="auto"/>).

$d = $word.ActiveDocument

foreach ($p in $d.Paragraphs) {
    $p.LineSpacingRule = 0   # wdLineSpaceSingle
}
